$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab: "Sheet1" -> "Hoja1" ---
$ws.Name = "Hoja1"

# --- Move the active selection on the sheet: C7 -> E6 ---
$ws.Range("E6").Select()

# --- Narrow/adjust column C's width slightly (closest width reachable
#     through the ColumnWidth object-model property) ---
$ws.Columns("C").ColumnWidth = 8

# --- Fix the stray password value in B2: "Patito00" -> "0" (kept as text,
#     not auto-converted to a number) ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0"
